$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-08 Saturday" "2025-02-09 Sunday"

Replace-Text "52×25=" "91×83="
Replace-Text "56×41=" "82×44="
Replace-Text "32×25=" "36×66="
Replace-Text "96×39=" "83×11="
Replace-Text "41×28=" "61×23="

Replace-Text "13×15=" "75×18="
Replace-Text "45×37=" "52×91="
Replace-Text "74×84=" "47×97="
Replace-Text "45×61=" "31×71="
Replace-Text "24×20=" "73×76="

Replace-Text "84×92=" "93×78="
Replace-Text "15×75=" "40×51="
Replace-Text "56×50=" "69×11="
Replace-Text "36×71=" "96×42="
Replace-Text "95×71=" "41×56="

Replace-Text "99×77=" "18×22="
Replace-Text "47×68=" "90×44="
Replace-Text "28×79=" "24×60="
Replace-Text "70×41=" "55×23="
Replace-Text "20×32=" "54×92="

Replace-Text "76×86=" "53×17="
Replace-Text "56×20=" "28×19="
Replace-Text "45×40=" "35×29="
Replace-Text "16×15=" "84×60="
Replace-Text "46×61=" "17×44="
